$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '63.696.17'
$ws.Range("E2").Value = '  -1.01%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.080.02'
$ws.Range("E3").Value = '  -0.95%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '593.70'
$ws.Range("E5").Value = '  +0.77%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '155.42'
$ws.Range("E6").Value = '  +2.26%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '3.079.17'
$ws.Range("E9").Value = '  -0.92%  '
$ws.Range("E10").Value = '  -1.04%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.94'
$ws.Range("E11").Value = '  +0.09%  '
$ws.Range("E12").Value = '  -1.71%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000238'
$ws.Range("E13").Value = '  -1.99%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '36.82'
$ws.Range("E14").Value = '  -2.81%  '
$ws.Range("E15").Value = '  +1.13%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.583.48'
$ws.Range("E16").Value = '  -1.25%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '7.20'
$ws.Range("E17").Value = '  -0.31%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '63.535.11'
$ws.Range("E18").Value = '  -0.62%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.074.81'
$ws.Range("E19").Value = '  -1.27%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '482.64'
$ws.Range("E20").Value = '  +2.77%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.51'
$ws.Range("E21").Value = '  -2.75%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.711'
$ws.Range("E22").Value = '  -3.64%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.57'
$ws.Range("E23").Value = '  -0.15%  '
$ws.Range("E24").Value = '  +1.11%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '81.77'
$ws.Range("E25").Value = '  +0.11%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '12.86'
$ws.Range("E26").Value = '  -3.16%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.80'
$ws.Range("E27").Value = '  +10.84%  '
$ws.Range("E28").Value = '  +0.17%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.69'
$ws.Range("E29").Value = '  +3.99%  '
$ws.Range("B30").Value = 'ImmutableX'
$ws.Range("C30").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.24'
$ws.Range("E30").Value = '  +1.74%  '
$ws.Range("B31").Value = 'PancakeSwap'
$ws.Range("C31").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.70'
$ws.Range("E31").Value = '  +0.02%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.999'
$ws.Range("E32").Value = '  -0.92%  '
$ws.Range("E33").Value = '  -3.19%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '27.27'
$ws.Range("E34").Value = '  -0.70%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0₃0831'
$ws.Range("E35").Value = '  -1.93%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.07'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '6.08'
$ws.Range("E37").Value = '  -1.12%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.30'
$ws.Range("E38").Value = '  -1.93%  '
$ws.Range("E39").Value = '  -1.60%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '9.24'
$ws.Range("E40").Value = '  -1.17%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '50.26'
$ws.Range("E41").Value = '  -1.19%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '445.49'
$ws.Range("E42").Value = '  -1.30%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.292'
$ws.Range("E43").Value = '  +0.53%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.113'
$ws.Range("E44").Value = '  +3.34%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0363'
$ws.Range("E45").Value = '  -1.89%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '40.01'
$ws.Range("E46").Value = '  +3.94%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.826.72'
$ws.Range("E47").Value = '  -0.79%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '132.53'
$ws.Range("E48").Value = '  +2.10%  '
$ws.Range("B49").Value = 'InjectiveProtocol'
$ws.Range("C49").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '25.32'
$ws.Range("E49").Value = '  +0.62%  '
$ws.Range("B50").Value = 'USDe'
$ws.Range("C50").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.999'
$ws.Range("E50").Value = '  +0.01%  '
$ws.Range("E51").Value = '  -1.32%  '
